$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the CORREL formula in D3 (correlation between Column1 and Column2 data ranges)
$ws.Range("D3").Formula = "=CORREL(A3:A12,B3:B12)"

# After typing a formula and pressing Enter, Excel moves the active cell down to D4
$ws.Range("D4").Select()
